$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 2 and row 4 for columns D, J, K, L, M, N, O, P, Q
$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr4 = "$col" + "4"
    $val2 = $ws.Range($addr2).Value2
    $val4 = $ws.Range($addr4).Value2
    $ws.Range($addr2).Value2 = $val4
    $ws.Range($addr4).Value2 = $val2
}
